# Update party abbreviation headers to include full party names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "CDU - Christian Democratic Union (Christlich Demokratische Union , CDU)"
$ws.Range("C1").Value = "CSU - Christian Social Union  (Christlich Soziale Union, CSU)"
$ws.Range("D1").Value = "FDP - Free Democrats  (Freie Demokratische Partei, FDP)"
$ws.Range("E1").Value = "G - Alliance 90-Greens (Bundnis 90-Die Grunen, G)"
$ws.Range("F1").Value = "PDS - Party of Democratic Socialism (Partei des Demokratischen Sozialismus, PDS)"
$ws.Range("G1").Value = "SPD - Social Democrats (Sozialdemokratische Partei Deutschlands, SPD)"
$ws.Range("H1").Value = "Linke - The Left (Die Linke, Linke)"
$ws.Range("I1").Value = "AfD - Alternative for Germany (Alternative für Deutschland, AfD)"
